$d = $word.ActiveDocument

# 1) Update the two representative names in paragraphs 3.1 and 3.2
$d.Content.Find.Execute("Javier Jiménez", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Betzabet Marín", 2)

$d.Content.Find.Execute("Fernando Gómez", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Araceli Becerril", 2)

# 2) Move the "_GoBack" bookmark from the end of the document (after
#    "subcontratación") to the start of the paragraph beginning
#    "4-. Debes de conocer la carta de VALORES ..."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute("4-. Debes de conocer", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$target.Collapse(1)

$d.Bookmarks.Add("_GoBack", $target)
